$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Preserve text formatting for Price/Volume columns so numeric-looking
# strings (e.g. "241.81") are not converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '36.554.27'
$ws.Cells.Item(2, 5).Value = '  -1.71%  '
$ws.Cells.Item(3, 4).Value = '2.062.37'
$ws.Cells.Item(3, 5).Value = '  +0.14%  '
$ws.Cells.Item(4, 5).Value = '  -0.03%  '
$ws.Cells.Item(5, 4).Value = '241.81'
$ws.Cells.Item(5, 5).Value = '  -2.90%  '
$ws.Cells.Item(6, 4).Value = '0.660'
$ws.Cells.Item(6, 5).Value = '  -1.10%  '
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 4).Value = '52.80'
$ws.Cells.Item(8, 5).Value = '  -7.48%  '
$ws.Cells.Item(9, 4).Value = '59.16'
$ws.Cells.Item(9, 5).Value = '  -1.91%  '
$ws.Cells.Item(10, 5).Value = '  -6.97%  '
$ws.Cells.Item(11, 4).Value = '0.0749'
$ws.Cells.Item(11, 5).Value = '  -4.80%  '
$ws.Cells.Item(12, 5).Value = '  -0.71%  '
$ws.Cells.Item(13, 4).Value = '0.915'
$ws.Cells.Item(13, 5).Value = '  -0.55%  '
$ws.Cells.Item(14, 4).Value = '14.68'
$ws.Cells.Item(14, 5).Value = '  -9.69%  '
$ws.Cells.Item(15, 4).Value = '2.360.52'
$ws.Cells.Item(15, 5).Value = '  -0.01%  '
$ws.Cells.Item(16, 4).Value = '5.40'
$ws.Cells.Item(16, 5).Value = '  -6.21%  '
$ws.Cells.Item(17, 4).Value = '2.074.42'
$ws.Cells.Item(17, 5).Value = '  +0.78%  '
$ws.Cells.Item(18, 4).Value = '36.473.68'
$ws.Cells.Item(18, 5).Value = '  -1.98%  '
$ws.Cells.Item(19, 4).Value = '16.34'
$ws.Cells.Item(19, 5).Value = '  -12.99%  '
$ws.Cells.Item(20, 4).Value = '71.61'
$ws.Cells.Item(20, 5).Value = '  -4.46%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0858'
$ws.Cells.Item(21, 5).Value = '  -4.96%  '
$ws.Cells.Item(22, 4).Value = '236.60'
$ws.Cells.Item(22, 5).Value = '  -0.85%  '
$ws.Cells.Item(23, 4).Value = '5.24'
$ws.Cells.Item(23, 5).Value = '  -4.45%  '
$ws.Cells.Item(24, 5).Value = '  +0.14%  '
$ws.Cells.Item(25, 5).Value = '  -4.97%  '
$ws.Cells.Item(26, 4).Value = '9.47'
$ws.Cells.Item(26, 5).Value = '  -2.08%  '
$ws.Cells.Item(27, 4).Value = '2.12'
$ws.Cells.Item(27, 5).Value = '  -3.28%  '
$ws.Cells.Item(28, 4).Value = '164.26'
$ws.Cells.Item(28, 5).Value = '  -3.63%  '
$ws.Cells.Item(29, 4).Value = '20.39'
$ws.Cells.Item(29, 5).Value = '  +0.71%  '
$ws.Cells.Item(30, 5).Value = '  -2.74%  '
$ws.Cells.Item(31, 4).Value = '5.05'
$ws.Cells.Item(31, 5).Value = '  -2.17%  '
$ws.Cells.Item(32, 4).Value = '1.14'
$ws.Cells.Item(32, 5).Value = '  -2.75%  '
$ws.Cells.Item(33, 5).Value = '  -1.58%  '
$ws.Cells.Item(34, 4).Value = '0.0594'
$ws.Cells.Item(34, 5).Value = '  -5.16%  '
$ws.Cells.Item(35, 5).Value = '  +0.09%  '
$ws.Cells.Item(36, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(36, 4).Value = '2.30'
$ws.Cells.Item(36, 5).Value = '  -0.43%  '
$ws.Cells.Item(37, 2).Value = 'WEMIXToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(37, 4).Value = '1.83'
$ws.Cells.Item(37, 5).Value = '  +3.37%  '
$ws.Cells.Item(38, 4).Value = '0.0821'
$ws.Cells.Item(38, 5).Value = '  -8.07%  '
$ws.Cells.Item(39, 5).Value = '  -7.87%  '
$ws.Cells.Item(40, 5).Value = '  -5.36%  '
$ws.Cells.Item(41, 4).Value = '4.82'
$ws.Cells.Item(41, 5).Value = '  -9.39%  '
$ws.Cells.Item(42, 5).Value = '  -2.98%  '
$ws.Cells.Item(43, 4).Value = '0.0216'
$ws.Cells.Item(43, 5).Value = '  -3.81%  '
$ws.Cells.Item(44, 4).Value = '0.0939'
$ws.Cells.Item(44, 5).Value = '  -6.85%  '
$ws.Cells.Item(45, 4).Value = '93.88'
$ws.Cells.Item(45, 5).Value = '  -3.11%  '
$ws.Cells.Item(46, 4).Value = '1.384.17'
$ws.Cells.Item(46, 5).Value = '  +8.20%  '
$ws.Cells.Item(47, 4).Value = '7.46'
$ws.Cells.Item(47, 5).Value = '  +8.71%  '
$ws.Cells.Item(48, 4).Value = '15.49'
$ws.Cells.Item(48, 5).Value = '  -11.99%  '
$ws.Cells.Item(49, 4).Value = '2.36'
$ws.Cells.Item(49, 5).Value = '  -3.38%  '
$ws.Cells.Item(50, 4).Value = '2.84'
$ws.Cells.Item(50, 5).Value = '  -0.71%  '
$ws.Cells.Item(51, 4).Value = '2.251.89'
$ws.Cells.Item(51, 5).Value = '  +0.03%  '
